$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 12
$ws.Range("C2").Value = 340
$ws.Range("C3").Value = 220
$ws.Range("B4").Value = 3600
$ws.Range("C4").Value = 0.001
$ws.Range("C5").Value = 125

$ws.Range("B5").Select()
